$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds numeric-looking text (e.g. "206.95").
# Excel auto-converts such strings to numbers on assignment, which loses
# the exact text formatting (trailing zeros, thousands-dot grouping, etc.)
# stored in the source data. Force those specific cells to Text format
# first so the literal string is preserved, matching the source sheet.

$ws.Range("D2").Value = '26.896.20'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.548.61'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.95'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.42'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0584'
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = '1.769.27'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").Value = '1.546.04'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("D16").Value = '26.908.17'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.57'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.75'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '0.0₃0685'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.01'
$ws.Range("E22").Value = '  -2.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.05'
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.89'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("D33").Value = '1.368.58'
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.95'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("E35").Value = '  +1.02%  '
$ws.Range("E36").Value = '  +3.85%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0164'
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  +6.94%  '
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("E44").Value = '  +2.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.63'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("B47").Value = 'mCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.28'
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.683.85'
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.63'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("E50").Value = '  +1.28%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0971'
$ws.Range("E51").Value = '  -1.57%  '
